# FamilyTree.xlsx fix: a row (ID=3, "Jonas") was missing from the sheet
# because of a CSV/Excel load bug. This restores it in correct ID order
# and patches up the parent/children/spouse cross references that
# referred to him.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 -- this pushes the existing rows 5-10
# (Edmond, Vere, Geneva, Blanche, Tammy, Jenny) down to rows 6-11 and
# keeps all their data intact.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with Jonas's record (ID 3), child of
# Alex (ID 2) and Vere (ID 5), spouse of Blanche (ID 7).
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Jonas"
$ws.Cells.Item(5, 3).Value = "Castro"
$ws.Cells.Item(5, 4).Value = "Utena"
$ws.Cells.Item(5, 5).Value = "2019-02-02"
$ws.Cells.Item(5, 6).Value = "2, 5"
$ws.Cells.Item(5, 7).Value = ""
$ws.Cells.Item(5, 8).Value = 7
$ws.Cells.Item(5, 9).Value = 2

# Alex (row 4) now also has Jonas as a child, alongside Geneva.
$ws.Cells.Item(4, 7).Value = "3, 6"

# Vere (now row 7 after the insert) likewise gains Jonas as a child.
$ws.Cells.Item(7, 7).Value = "3, 6"

# Blanche (now row 9) is Jonas's spouse.
$ws.Cells.Item(9, 8).Value = 3

# Tammy (now row 10) has both Edmond and Jenny listed as parents.
$ws.Cells.Item(10, 6).Value = "4, 9"
